$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.778.59'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.454.88'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +6.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '480.24'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +14.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +7.29%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.447.72'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0968'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +11.35%  '
$ws.Range('E11').Value = '  +3.42%  '
$ws.Range('E12').Value = '  +6.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.123'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.845.94'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '54.896.22'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.44'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E17').Value = '  +14.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.442.32'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.32'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +9.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '313.99'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +5.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.87'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +10.92%  '
$ws.Range('E22').Value = '  -0.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '57.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.164'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +14.52%  '
$ws.Range('E26').Value = '  +10.62%  '
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.543.76'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.45%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.31'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.22%  '
$ws.Range('E30').Value = '  +17.71%  '
$ws.Range('E31').Value = '  -0.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.10'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '17.92'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.68%  '
$ws.Range('E34').Value = '  +10.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.12'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.07%  '
$ws.Range('E36').Value = '  +12.33%  '
$ws.Range('E37').Value = '  +6.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.846'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.91%  '
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '33.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.40'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.61%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0544'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.58%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.596'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.73%  '
$ws.Range('E44').Value = '  +10.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.63'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +10.78%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '254.44'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +26.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0902'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +10.14%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0222'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +8.62%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.929.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.01'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +9.37%  '
